$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.921.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.00%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.637.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.27%  '
# Row 4
$ws.Range("E4").Value = '  -0.03%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.01%  '
# Row 6
$ws.Range("E6").Value = '  -0.49%  '
# Row 7
$ws.Range("E7").Value = '  -0.03%  '
# Row 8
$ws.Range("E8").Value = '  +0.85%  '
# Row 9
$ws.Range("E9").Value = '  -0.27%  '
# Row 10
$ws.Range("E10").Value = '  -0.29%  '
# Row 11
$ws.Range("E11").Value = '  +0.30%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.869.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.26%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.642.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.61%  '
# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.565'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.13%  '
# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.89%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.20%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.932.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.04%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.45'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.07%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.24%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.07%  '
# Row 22
$ws.Range("E22").Value = '  +0.15%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.76%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.26%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.51%  '
# Row 26
$ws.Range("E26").Value = '  +0.16%  '
# Row 27
$ws.Range("E27").Value = '  +0.09%  '
# Row 28
$ws.Range("E28").Value = '  -0.47%  '
# Row 29
$ws.Range("E29").Value = '  +0.06%  '
# Row 30
$ws.Range("E30").Value = '  +0.14%  '
# Row 31
$ws.Range("E31").Value = '  -0.19%  '
# Row 32
$ws.Range("E32").Value = '  +1.17%  '
# Row 33
$ws.Range("E33").Value = '  +1.29%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.405.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.25%  '
# Row 35
$ws.Range("E35").Value = '  +3.21%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.04%  '
# Row 37
$ws.Range("E37").Value = '  -0.82%  '
# Row 38
$ws.Range("E38").Value = '  +0.75%  '
# Row 39
$ws.Range("E39").Value = '  -0.54%  '
# Row 40
$ws.Range("E40").Value = '  -1.92%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.02%  '
# Row 42
$ws.Range("E42").Value = '  -1.46%  '
# Row 43
$ws.Range("E43").Value = '  +2.21%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.99%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.34%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.777.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.12%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.54%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.75'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.17%  '
# Row 49
$ws.Range("E49").Value = '  +2.59%  '
# Row 50
$ws.Range("E50").Value = '  -0.26%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.81%  '
